$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.8
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.6
$ws.Range("L5").Value = 6
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 1.5
$ws.Range("U5").Value = 4.3
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 1.17
$ws.Range("AA5").Value = 2.38
$ws.Range("AB5").Value = 1.53
$ws.Range("AD5").Value = 7
$ws.Range("AF5").Value = 13
$ws.Range("AK5").Value = 23
$ws.Range("AL5").Value = 101
$ws.Range("AM5").Value = 9.5
$ws.Range("AN5").Value = 23
$ws.Range("AQ5").Value = 51

# Row 6
$ws.Range("G6").Value = 2.75
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 3.75
$ws.Range("L6").Value = 4
$ws.Range("N6").Value = 5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.54
$ws.Range("S6").Value = 3.5
$ws.Range("T6").Value = 1.3
$ws.Range("U6").Value = 6.2
$ws.Range("V6").Value = 1.13
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 1.08
$ws.Range("Y6").Value = 1.78
$ws.Range("Z6").Value = 2.03
$ws.Range("AG6").Value = 34
$ws.Range("AM6").Value = 6
$ws.Range("AQ6").Value = 34

# Row 9
$ws.Range("G9").Value = 3.4
$ws.Range("H9").Value = 2.9
$ws.Range("I9").Value = 2.3
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 1.95
$ws.Range("L9").Value = 3.2
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 2.63
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98
$ws.Range("S9").Value = 2.5
$ws.Range("T9").Value = 1.5
$ws.Range("AF9").Value = 34
$ws.Range("AK9").Value = 17
$ws.Range("AM9").Value = 6.5
$ws.Range("AN9").Value = 10

# Row 14
$ws.Range("G14").Value = 6.7
$ws.Range("H14").Value = 4.25
$ws.Range("I14").Value = 1.42
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 2.32
$ws.Range("L14").Value = 1.91
$ws.Range("P14").Value = 3.7
$ws.Range("S14").Value = 1.57
$ws.Range("T14").Value = 2.12
$ws.Range("W14").Value = 2.35
$ws.Range("X14").Value = 1.47
$ws.Range("AA14").Value = 1.75
$ws.Range("AB14").Value = 1.87
$ws.Range("AD14").Value = 50
$ws.Range("AE14").Value = 21
$ws.Range("AF14").Value = 150
$ws.Range("AG14").Value = 70
$ws.Range("AH14").Value = 55
$ws.Range("AI14").Value = 13.5
$ws.Range("AJ14").Value = 8.75
$ws.Range("AK14").Value = 16.5
$ws.Range("AL14").Value = 70
$ws.Range("AM14").Value = 7.7
$ws.Range("AN14").Value = 7.2
$ws.Range("AO14").Value = 8
$ws.Range("AP14").Value = 9.75
$ws.Range("AR14").Value = 23
$ws.Range("AS14").Value = 450

# Row 15
$ws.Range("G15").Value = 4.35
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 1.65
$ws.Range("J15").Value = 4.35
$ws.Range("K15").Value = 2.45
$ws.Range("L15").Value = 2.07
$ws.Range("P15").Value = 3.65
$ws.Range("S15").Value = 1.6
$ws.Range("T15").Value = 2.07
$ws.Range("W15").Value = 2.42
$ws.Range("X15").Value = 1.44
$ws.Range("AA15").Value = 1.65
$ws.Range("AB15").Value = 2
$ws.Range("AC15").Value = 14.5
$ws.Range("AD15").Value = 26
$ws.Range("AE15").Value = 14.5
$ws.Range("AF15").Value = 70
$ws.Range("AH15").Value = 40
$ws.Range("AI15").Value = 14
$ws.Range("AJ15").Value = 8
$ws.Range("AK15").Value = 15
$ws.Range("AL15").Value = 60
$ws.Range("AN15").Value = 8.5
$ws.Range("AO15").Value = 8.25
$ws.Range("AP15").Value = 12.5
$ws.Range("AQ15").Value = 12
$ws.Range("AR15").Value = 22
$ws.Range("AS15").Value = 400
